# TC-09 updated: remove the "project" column (column A) from the tc009
# test-data sheet, shifting epic/feature/id/reqId left by one column, and
# make tc009 the active sheet/tab.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("tc009")

# Delete the entire first column (the "project" / "STG- SPARK Modernization"
# column), shifting B:E left to A:D.
$ws.Columns.Item(1).Delete()

# Select/activate the now-4-column data, mirroring Excel's recorded
# selection after the edit.
$ws.Activate()
$ws.Range("B2").Select()

$wb.Save()
